# "Fruta / hortaliza, semanal"
#
# This adds a new weekly price-report entry for Cebollín (Vega Monumental
# Concepción) by inserting two new data rows right after the existing
# "Región de Ñuble" rows 54/55 (which carry the reusable 6-unit-package
# price pattern), shifting all subsequent rows down by two, and then
# updating the date (Fecha) and origin (Origen) of the two newly inserted
# rows to reflect the new report date / region.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 54 - everything from the
# old row 54 onward shifts down to make room (old 54..70 -> new 56..72).
$ws.Rows("54:55").Insert()

# The rows that used to be 54/55 are now at 56/57; duplicate them into the
# freshly inserted 54/55 slots so they start out with identical data.
$ws.Rows(56).Copy()
$ws.Rows(54).PasteSpecial()

$ws.Rows(57).Copy()
$ws.Rows(55).PasteSpecial()

# Update the new rows' Fecha (column D) and Origen (column O) to reflect
# the new weekly report entry.
$ws.Cells.Item(54, 4).Value2 = 44825
$ws.Cells.Item(54, 15).Value2 = "Región Metropolitana"

$ws.Cells.Item(55, 4).Value2 = 44825
$ws.Cells.Item(55, 15).Value2 = "Región Metropolitana"
